# Applies the "Auto report writing" fix described by the commit diff.
# Rows 2 & 3 (login success cases): Expected results / Actual results -> "Login Successfully"
# Rows 4-22 (login fail cases): Expected results / Actual results -> "Login Fail", and the
#   "Pass / Fail" column is corrected from "Fail" to "Pass" (the actual vs expected now match).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 3: successful login
foreach ($r in 2..3) {
    $ws.Cells.Item($r, 6).Value = "Login Successfully"   # F: Expected results
    $ws.Cells.Item($r, 7).Value = "Login Successfully"   # G: Actual results
    # H (Pass / Fail) stays "Pass" - unchanged
}

# Rows 4 through 22: failed login
foreach ($r in 4..22) {
    $ws.Cells.Item($r, 6).Value = "Login Fail"           # F: Expected results
    $ws.Cells.Item($r, 7).Value = "Login Fail"           # G: Actual results
    $ws.Cells.Item($r, 8).Value = "Pass"                 # H: Pass / Fail corrected to Pass
}
